# Edit script: applies the two logical changes described by the diff:
# 1. Fix a trailing-space typo "characters(Tom )" -> "characters(Tom)" (and the
#    associated ".1" duplicate-suffixed header) wherever it appears across the
#    workbook's worksheets.
# 2. Rebuild the "from_tom_and_jerry" sheet so that, instead of being pivoted
#    against the boolean "is_cat" feature, it is pivoted against every value of
#    the "characters" feature (one column per character), including fixing the
#    same "Tom " -> "Tom" typo in the new header row.

$wb = $excel.ActiveWorkbook

# --- 1. Global typo fix: "characters(Tom )" -> "characters(Tom)" -----------
foreach ($sheet in $wb.Worksheets) {
    [void]$sheet.Cells.Replace("characters(Tom )", "characters(Tom)")
}

# --- 2. Rebuild the from_tom_and_jerry sheet --------------------------------
$ws = $wb.Worksheets.Item("from_tom_and_jerry")

$headers = @("characters","characters(Agatsuma Zenitsu)","characters(Akai Shuuichi)","characters(Android 18)","characters(Ban)","characters(Beerus)","characters(Big Mom)","characters(Blue-Eyes White Dragon)","characters(Boros)","characters(Bulma)","characters(Cell)","characters(Champa)","characters(Daki)","characters(Dark Magician)","characters(Dark Magician Girl)","characters(Derieri)","characters(Diane)","characters(Doraemon)","characters(Edogawa Conan)","characters(Elizabeth Liones)","characters(Eren Yeager)","characters(Escanor)","characters(Frieza)","characters(Fubuki)","characters(Fushiguro Megumi)","characters(Galand)","characters(Garou)","characters(Genos)","characters(Gojo Satoru)","characters(Gowther)","characters(Haruno Sakura)","characters(Hashibira Inosuke)","characters(Hattori Heiji)","characters(Hawk)","characters(Howl Jenkins Pendragon)","characters(Hyuga Hinata)","characters(Itadori Yuji)","characters(Jerry)","characters(Kaiba Seto)","characters(Kaito Kid)","characters(Kamado Nezuko)","characters(Kamado Tanjiro)","characters(Kanroji Mitsuri)","characters(Katsuya Jonouchi)","characters(King)","characters(Kirito)","characters(Kocho Shinobu)","characters(Koro Sensei)","characters(Kudo Sinichi)","characters(Levi Ackerman)","characters(Majin Buu)","characters(Marik Ishtar)","characters(Melascula)","characters(Meliodas)","characters(Merlin)","characters(Mickey)","characters(Monkey D.Luffy)","characters(Muto Yugi)","characters(Nami)","characters(Nobita)","characters(Obelisk)","characters(Oggy)","characters(Osiris)","characters(Pegasus)","characters(Piccolo)","characters(Pikachu)","characters(Pink Panther)","characters(Po)","characters(Ra)","characters(Ran Mori)","characters(Rimuru Tempest)","characters(Ryomen Sukuna)","characters(Saitama)","characters(Sanji)","characters(Shanks)","characters(Shenron)","characters(Simba)","characters(Son Goku)","characters(Tatsumaki)","characters(Tokitou Muichirou)","characters(Tom)","characters(Tomioka Giyuu)","characters(Totoro)","characters(Uchiha Sasuke)","characters(Uzumaki Naruto)","characters(Vados)","characters(Vegeta)","characters(Whis)","characters(Yugi)","characters(Zeldris)","characters(Zeno)","characters(Zoro)")
$row2 = @(0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9910071942446044,0.9739583333333334,0.9866310160427808,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9866310160427808,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9866310160427808,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.0260416666666666,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9866310160427808,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9866310160427808,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9866310160427808,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9910071942446044,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9866310160427808,0.9739583333333334,0.9739583333333334,0.0260416666666666,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9866310160427808,0.9739583333333334,0.9739583333333334,0.9739583333333334,0.9866310160427808,0.9739583333333334)
$row3 = @(0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0089928057553956,0.0260416666666666,0.0133689839572192,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0133689839572192,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0133689839572192,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.9739583333333334,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0133689839572192,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0133689839572192,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0133689839572192,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0089928057553956,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0133689839572192,0.0260416666666666,0.0260416666666666,0.9739583333333334,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0133689839572192,0.0260416666666666,0.0260416666666666,0.0260416666666666,0.0133689839572192,0.0260416666666666)

$colCount = $headers.Length

# Row 1: headers ("characters", "characters(Agatsuma Zenitsu)", ...)
$headerArr = New-Object 'object[,]' 1,$colCount
for ($i = 0; $i -lt $colCount; $i++) {
    $headerArr[0,$i] = $headers[$i]
}
$headerRange = $ws.Range($ws.Cells.Item(1,1), $ws.Cells.Item(1,$colCount))
$headerRange.Value = $headerArr

# Make sure the newly-created header cells (D1 onward) pick up the same
# bold/centered/bordered style that A1:C1 already had.
[void]$ws.Cells.Item(1,1).Copy()
[void]$headerRange.PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Row 2 values (B2:CN2) -- A2 keeps its existing "from_tom_and_jerry(0)" label
$dataCount = $row2.Length
$row2Arr = New-Object 'object[,]' 1,$dataCount
for ($i = 0; $i -lt $dataCount; $i++) {
    $row2Arr[0,$i] = $row2[$i]
}
$row2Range = $ws.Range($ws.Cells.Item(2,2), $ws.Cells.Item(2,$colCount))
$row2Range.Value = $row2Arr

# Row 3 values (B3:CN3) -- A3 keeps its existing "from_tom_and_jerry(1)" label
$row3Arr = New-Object 'object[,]' 1,$dataCount
for ($i = 0; $i -lt $dataCount; $i++) {
    $row3Arr[0,$i] = $row3[$i]
}
$row3Range = $ws.Range($ws.Cells.Item(3,2), $ws.Cells.Item(3,$colCount))
$row3Range.Value = $row3Arr
